{"js": "// Update the date line and every \"a\u00d7b=c\" answer cell to the new values\n// described by the commit (output regenerated for 2025-06-30 Monday).\nconst replacements = [\n  [\"2025-06-29 Sunday\", \"2025-06-30 Monday\"],\n  [\"711\u00d74=2844\", \"513\u00d78=4104\"],\n  [\"461\u00d78=3688\", \"541\u00d73=1623\"],\n  [\"495\u00d76=2970\", \"308\u00d72=616\"],\n  [\"493\u00d73=1479\", \"765\u00d79=6885\"],\n  [\"284\u00d79=2556\", \"539\u00d72=1078\"],\n  [\"999\u00d74=3996\", \"686\u00d74=2744\"],\n  [\"746\u00d75=3730\", \"430\u00d74=1720\"],\n  [\"364\u00d73=1092\", \"696\u00d79=6264\"],\n  [\"609\u00d72=1218\", \"229\u00d74=916\"],\n  [\"120\u00d74=480\", \"674\u00d77=4718\"],\n  [\"577\u00d74=2308\", \"397\u00d76=2382\"],\n  [\"113\u00d78=904\", \"913\u00d78=7304\"],\n  [\"767\u00d77=5369\", \"453\u00d76=2718\"],\n  [\"166\u00d76=996\", \"358\u00d72=716\"],\n  [\"470\u00d74=1880\", \"426\u00d73=1278\"],\n  [\"630\u00d76=3780\", \"738\u00d77=5166\"],\n  [\"511\u00d74=2044\", \"712\u00d76=4272\"],\n  [\"154\u00d76=924\", \"813\u00d75=4065\"],\n  [\"816\u00d72=1632\", \"391\u00d75=1955\"],\n  [\"517\u00d77=3619\", \"466\u00d76=2796\"],\n  [\"327\u00d77=2289\", \"961\u00d73=2883\"],\n  [\"937\u00d78=7496\", \"929\u00d74=3716\"],\n  [\"402\u00d73=1206\", \"206\u00d78=1648\"],\n  [\"552\u00d78=4416\", \"186\u00d78=1488\"],\n  [\"679\u00d76=4074\", \"766\u00d77=5362\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and every \"a\u00d7b=c\" answer cell to the new values\n# described by the commit (output regenerated for 2025-06-30 Monday).\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-06-29 Sunday\", \"2025-06-30 Monday\"),\n    @(\"711\u00d74=2844\", \"513\u00d78=4104\"),\n    @(\"461\u00d78=3688\", \"541\u00d73=1623\"),\n    @(\"495\u00d76=2970\", \"308\u00d72=616\"),\n    @(\"493\u00d73=1479\", \"765\u00d79=6885\"),\n    @(\"284\u00d79=2556\", \"539\u00d72=1078\"),\n    @(\"999\u00d74=3996\", \"686\u00d74=2744\"),\n    @(\"746\u00d75=3730\", \"430\u00d74=1720\"),\n    @(\"364\u00d73=1092\", \"696\u00d79=6264\"),\n    @(\"609\u00d72=1218\", \"229\u00d74=916\"),\n    @(\"120\u00d74=480\", \"674\u00d77=4718\"),\n    @(\"577\u00d74=2308\", \"397\u00d76=2382\"),\n    @(\"113\u00d78=904\", \"913\u00d78=7304\"),\n    @(\"767\u00d77=5369\", \"453\u00d76=2718\"),\n    @(\"166\u00d76=996\", \"358\u00d72=716\"),\n    @(\"470\u00d74=1880\", \"426\u00d73=1278\"),\n    @(\"630\u00d76=3780\", \"738\u00d77=5166\"),\n    @(\"511\u00d74=2044\", \"712\u00d76=4272\"),\n    @(\"154\u00d76=924\", \"813\u00d75=4065\"),\n    @(\"816\u00d72=1632\", \"391\u00d75=1955\"),\n    @(\"517\u00d77=3619\", \"466\u00d76=2796\"),\n    @(\"327\u00d77=2289\", \"961\u00d73=2883\"),\n    @(\"937\u00d78=7496\", \"929\u00d74=3716\"),\n    @(\"402\u00d73=1206\", \"206\u00d78=1648\"),\n    @(\"552\u00d78=4416\", \"186\u00d78=1488\"),\n    @(\"679\u00d76=4074\", \"766\u00d77=5362\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
